$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "nft1.jpg"
$ws.Range("B2").Value = "https://cloudflare-ipfs.com/ipfs/QmV4WJFR7ckRJrJqpvR4bZzNjZNTbV5DWhfXbtyTqRkX7e"

$ws.Range("A3").Value = "nft2.jpg"
$ws.Range("B3").Value = "https://cloudflare-ipfs.com/ipfs/QmSpe6oTxU58utX2aoLXZg1e2Q3LsR6Hygv5b6CTPRCMqR"

$ws.Range("A4").Value = "nft3.jpg"
$ws.Range("B4").Value = "https://cloudflare-ipfs.com/ipfs/QmVoRVn2NbX8gQg6LbjKtDc9NVbcjTa6WoTQxqcyQuw9DB"
